# PlayerDB_Sheet.xlsx : add a "description" column after "prefabPath"
# (commit: "Fix : GameOver, player이동 수정" / insert description data for
# the two existing Player rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (prefabPath is C, job used to be D) and give it
# the same width as the column it was split off from (C).
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Header + the two data rows for the new "description" column.
$ws.Range("D1").Value = "description"
$ws.Range("D2").Value = "단단하다"
$ws.Range("D3").Value = "날렵하다"

# Restore the selection to where the author left off.
$ws.Range("F10").Select() | Out-Null
